# Apply updated "dSF" (column F) values pulled from the repulled source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = -4
    4  = -1
    5  = 1
    6  = 2
    7  = 2
    8  = 4
    9  = 3
    10 = 2
    11 = 4
    12 = 2
    14 = -4
    15 = 2
    16 = 3
    17 = 1
    19 = -1
    20 = -2
    22 = -3
    23 = -6
    24 = 3
    25 = -4
    26 = 2
    28 = 3
    29 = 1
    30 = -2
    31 = 2
    32 = 10
    33 = 4
    35 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
